# Update cryptos list cell values per upstream data refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.163.81"
$ws.Range("E2").Value = "  -1.80%  "
$ws.Range("D3").Value = "1.823.11"
$ws.Range("E3").Value = "  -1.27%  "
$ws.Range("E4").Value = "  -0.50%  "
$ws.Range("D5").Value = "'312.84"
$ws.Range("E5").Value = "  -1.82%  "
$ws.Range("E6").Value = "  -0.47%  "
$ws.Range("E7").Value = "  -1.78%  "
$ws.Range("D8").Value = "'0.3696"
$ws.Range("E8").Value = "  -1.41%  "
$ws.Range("D9").Value = "'0.07247"
$ws.Range("E9").Value = "  -1.35%  "
$ws.Range("D10").Value = "'0.8583"
$ws.Range("E10").Value = "  -2.45%  "
$ws.Range("D11").Value = "'20.97"
$ws.Range("E11").Value = "  -2.82%  "
$ws.Range("D12").Value = "1.816.89"
$ws.Range("E12").Value = "  -1.67%  "
$ws.Range("D13").Value = "'6.711"
$ws.Range("E13").Value = "  -0.28%  "
$ws.Range("D14").Value = "'0.07092"
$ws.Range("E14").Value = "  -0.10%  "
$ws.Range("E15").Value = "  -2.70%  "
$ws.Range("D16").Value = "'89.57"
$ws.Range("E16").Value = "  +2.18%  "
$ws.Range("E17").Value = "  -0.71%  "
$ws.Range("D18").Value = "'0.000008852"
$ws.Range("E18").Value = "  -1.45%  "
$ws.Range("E19").Value = "  -0.52%  "
$ws.Range("B20").Value = "BitDAO"
$ws.Range("C20").Value = "https://coinranking.com/coin/N2IgQ9Xme+bitdao-bit"
$ws.Range("D20").Value = "'0.5111"
$ws.Range("E20").Value = "  +2.50%  "
$ws.Range("B21").Value = "Avalanche"
$ws.Range("C21").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D21").Value = "'15.04"
$ws.Range("E21").Value = "  -2.64%  "
$ws.Range("B22").Value = "WrappedBTC"
$ws.Range("C22").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D22").Value = "27.244.89"
$ws.Range("E22").Value = "  -1.56%  "
$ws.Range("B23").Value = "Uniswap"
$ws.Range("C23").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D23").Value = "'5.131"
$ws.Range("E23").Value = "  -2.38%  "
$ws.Range("B24").Value = "Cosmos"
$ws.Range("C24").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D24").Value = "'10.91"
$ws.Range("E24").Value = "  -2.27%  "
$ws.Range("B25").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C25").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D25").Value = "2.054.34"
$ws.Range("E25").Value = "  -1.29%  "
$ws.Range("B26").Value = "Toncoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D26").Value = "'1.988"
$ws.Range("E26").Value = "  -2.41%  "
$ws.Range("B27").Value = "Monero"
$ws.Range("C27").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D27").Value = "'152.46"
$ws.Range("E27").Value = "  -1.96%  "
$ws.Range("B28").Value = "LidoDAOToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D28").Value = "'2.189"
$ws.Range("E28").Value = "  +2.36%  "
$ws.Range("B29").Value = "EthereumClassic"
$ws.Range("C29").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D29").Value = "'18.40"
$ws.Range("E29").Value = "  -0.75%  "
$ws.Range("B30").Value = "InternetComputer(DFINITY)"
$ws.Range("C30").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D30").Value = "'5.224"
$ws.Range("E30").Value = "  -2.92%  "
$ws.Range("B31").Value = "BitcoinCash"
$ws.Range("C31").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D31").Value = "'116.38"
$ws.Range("E31").Value = "  -3.31%  "
$ws.Range("B32").Value = "Stellar"
$ws.Range("C32").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D32").Value = "'0.08841"
$ws.Range("E32").Value = "  -0.93%  "
$ws.Range("B33").Value = "ARBITRUM"
$ws.Range("C33").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D33").Value = "'1.189"
$ws.Range("E33").Value = "  -3.14%  "
$ws.Range("B34").Value = "ImmutableX"
$ws.Range("C34").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D34").Value = "'0.7498"
$ws.Range("E34").Value = "  -3.46%  "
$ws.Range("B35").Value = "Filecoin"
$ws.Range("C35").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D35").Value = "'4.437"
$ws.Range("E35").Value = "  -2.57%  "
$ws.Range("B36").Value = "HuobiToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D36").Value = "'2.805"
$ws.Range("E36").Value = "  -3.52%  "
$ws.Range("B37").Value = "Frax"
$ws.Range("C37").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D37").Value = "'1.005"
$ws.Range("E37").Value = "  -0.51%  "
$ws.Range("B38").Value = "TrustWalletToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D38").Value = "'1.115"
$ws.Range("E38").Value = "  -2.24%  "
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").Value = "'0.01969"
$ws.Range("E39").Value = "  -0.04%  "
$ws.Range("B40").Value = "Hedera"
$ws.Range("C40").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D40").Value = "'0.05249"
$ws.Range("E40").Value = "  -1.68%  "
$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").Value = "'7.301"
$ws.Range("E41").Value = "  +1.45%  "
$ws.Range("B42").Value = "MXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D42").Value = "'2.877"
$ws.Range("E42").Value = "  +0.47%  "
$ws.Range("B43").Value = "Algorand"
$ws.Range("C43").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D43").Value = "'0.1692"
$ws.Range("E43").Value = "  +0.90%  "
$ws.Range("B44").Value = "TheSandbox"
$ws.Range("C44").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D44").Value = "'0.5034"
$ws.Range("E44").Value = "  -2.48%  "
$ws.Range("B45").Value = "Aptos"
$ws.Range("C45").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D45").Value = "'8.665"
$ws.Range("E45").Value = "  -2.98%  "
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").Value = "'10.59"
$ws.Range("E46").Value = "  -1.10%  "
$ws.Range("B47").Value = "Quant"
$ws.Range("C47").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D47").Value = "'106.50"
$ws.Range("E47").Value = "  -3.65%  "
$ws.Range("B48").Value = "Decentraland"
$ws.Range("C48").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D48").Value = "'0.4739"
$ws.Range("E48").Value = "  +0.22%  "
$ws.Range("B49").Value = "PaxDollar"
$ws.Range("C49").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D49").Value = "'1.005"
$ws.Range("E49").Value = "  -0.48%  "
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").Value = "'0.06390"
$ws.Range("E50").Value = "  -1.62%  "
$ws.Range("B51").Value = "NEARProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D51").Value = "'1.667"
$ws.Range("E51").Value = "  -1.90%  "
